$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Database" column header in J1
$ws.Range("J1").Value = "Database"

# Fill J2:J223 with the database name "Econ2011"
$ws.Range("J2:J223").Value = "Econ2011"
